# Update workbook with carjacking data through 2022-09-14 (add data for 2022-09-22 run)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-09-14"

# Update the September label in column A, row 10
$ws.Range("A10").Value = "September (through 09-14)"

# Update September (row 10) figures for years 2015-2022 (columns B-I)
$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 24
$ws.Range("D10").Value = 35
$ws.Range("E10").Value = 27
$ws.Range("F10").Value = 30
$ws.Range("G10").Value = 53
$ws.Range("H10").Value = 72
$ws.Range("I10").Value = 65

# Update Total (row 11) figures for years 2015-2022 (columns B-I)
$ws.Range("B11").Value = 209
$ws.Range("C11").Value = 405
$ws.Range("D11").Value = 586
$ws.Range("E11").Value = 517
$ws.Range("F11").Value = 379
$ws.Range("G11").Value = 837
$ws.Range("H11").Value = 1142
$ws.Range("I11").Value = 1202
